# Update "想去人数" (want-to-go count) figures on the sheets that list
# individual exhibitions: "展览" and "全部类型".
#
#   Row 5  (曙光次元动漫游戏嘉年华)            F: 39   -> 41
#   Row 10 (第十五届次元之门动漫游戏博览会)      F: 5633 -> 5641
#   Row 11 (首届AT次元时代动漫游戏嘉年华)        F: 5001 -> 5007

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F5").Value = 41
    $ws.Range("F10").Value = 5641
    $ws.Range("F11").Value = 5007
}
